$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new data rows before row 260, shifting existing rows 260-327 down to 263-330.
$ws.Rows("260:262").Insert()

# Common (unchanged) field values shared by every data row in this block.
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$tipo      = "Fruta"
$prodId    = 100101
$producto  = "Berries"
$catId     = 100101007
$categoria = "Kiwi"
$origen    = "Región de O'Higgins"

# New row 260
$ws.Cells.Item(260, 1).Value  = $mercadoId
$ws.Cells.Item(260, 2).Value  = $mercado
$ws.Cells.Item(260, 3).Value  = $region
$ws.Cells.Item(260, 4).Value  = 44468
$ws.Cells.Item(260, 5).Value  = $codreg
$ws.Cells.Item(260, 6).Value  = $tipo
$ws.Cells.Item(260, 7).Value  = $prodId
$ws.Cells.Item(260, 8).Value  = $producto
$ws.Cells.Item(260, 9).Value  = $catId
$ws.Cells.Item(260, 10).Value = $categoria
$ws.Cells.Item(260, 11).Value = "Hayward"
$ws.Cells.Item(260, 12).Value = "Especial"
$ws.Cells.Item(260, 13).Value = 50
$ws.Cells.Item(260, 14).Value = 20000
$ws.Cells.Item(260, 15).Value = 20000
$ws.Cells.Item(260, 16).Value = 20000
$ws.Cells.Item(260, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(260, 18).Value = $origen
$ws.Cells.Item(260, 19).Value = 1333
$ws.Cells.Item(260, 20).Value = 15

# New row 261
$ws.Cells.Item(261, 1).Value  = $mercadoId
$ws.Cells.Item(261, 2).Value  = $mercado
$ws.Cells.Item(261, 3).Value  = $region
$ws.Cells.Item(261, 4).Value  = 44468
$ws.Cells.Item(261, 5).Value  = $codreg
$ws.Cells.Item(261, 6).Value  = $tipo
$ws.Cells.Item(261, 7).Value  = $prodId
$ws.Cells.Item(261, 8).Value  = $producto
$ws.Cells.Item(261, 9).Value  = $catId
$ws.Cells.Item(261, 10).Value = $categoria
$ws.Cells.Item(261, 11).Value = "Hayward"
$ws.Cells.Item(261, 12).Value = "Primera"
$ws.Cells.Item(261, 13).Value = 40
$ws.Cells.Item(261, 14).Value = 10000
$ws.Cells.Item(261, 15).Value = 10000
$ws.Cells.Item(261, 16).Value = 10000
$ws.Cells.Item(261, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(261, 18).Value = $origen
$ws.Cells.Item(261, 19).Value = 1000
$ws.Cells.Item(261, 20).Value = 10

# New row 262
$ws.Cells.Item(262, 1).Value  = $mercadoId
$ws.Cells.Item(262, 2).Value  = $mercado
$ws.Cells.Item(262, 3).Value  = $region
$ws.Cells.Item(262, 4).Value  = 44468
$ws.Cells.Item(262, 5).Value  = $codreg
$ws.Cells.Item(262, 6).Value  = $tipo
$ws.Cells.Item(262, 7).Value  = $prodId
$ws.Cells.Item(262, 8).Value  = $producto
$ws.Cells.Item(262, 9).Value  = $catId
$ws.Cells.Item(262, 10).Value = $categoria
$ws.Cells.Item(262, 11).Value = "Hayward"
$ws.Cells.Item(262, 12).Value = "Primera"
$ws.Cells.Item(262, 13).Value = 40
$ws.Cells.Item(262, 14).Value = 22000
$ws.Cells.Item(262, 15).Value = 22000
$ws.Cells.Item(262, 16).Value = 22000
$ws.Cells.Item(262, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(262, 18).Value = $origen
$ws.Cells.Item(262, 19).Value = 1222
$ws.Cells.Item(262, 20).Value = 18

Write-Host "Done"
